$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the Price column as Text so numeric-looking values
# ("1.002", "5.818", etc.) are preserved as literal strings rather
# than being coerced into numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.588.91"
$ws.Range("E2").Value = "  -3.37%  "

$ws.Range("D3").Value = "1.848.89"
$ws.Range("E3").Value = "  -3.90%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -1.12%  "

$ws.Range("D5").Value = "335.38"
$ws.Range("E5").Value = "  +2.70%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -1.00%  "

$ws.Range("D7").Value = "0.4653"
$ws.Range("E7").Value = "  -3.47%  "

$ws.Range("D8").Value = "0.3902"
$ws.Range("E8").Value = "  -3.68%  "

$ws.Range("D9").Value = "46.19"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("D10").Value = "0.07894"
$ws.Range("E10").Value = "  -3.57%  "

$ws.Range("E11").Value = "  -3.19%  "

$ws.Range("D12").Value = "22.24"
$ws.Range("E12").Value = "  -6.09%  "

$ws.Range("D13").Value = "1.825.58"
$ws.Range("E13").Value = "  -4.96%  "

$ws.Range("D14").Value = "5.818"
$ws.Range("E14").Value = "  -3.94%  "

$ws.Range("D15").Value = "6.961"
$ws.Range("E15").Value = "  -4.33%  "

$ws.Range("D16").Value = "0.06899"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "87.87"
$ws.Range("E17").Value = "  -3.80%  "

$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("E19").Value = "  -3.47%  "

$ws.Range("D20").Value = "17.05"
$ws.Range("E20").Value = "  -3.03%  "

$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "28.607.76"
$ws.Range("E22").Value = "  -3.29%  "

$ws.Range("D23").Value = "5.383"
$ws.Range("E23").Value = "  -4.52%  "

$ws.Range("E24").Value = "  -6.22%  "

$ws.Range("D25").Value = "2.179"
$ws.Range("E25").Value = "  -0.99%  "

$ws.Range("D26").Value = "2.075.73"
$ws.Range("E26").Value = "  -3.48%  "

$ws.Range("D27").Value = "153.29"
$ws.Range("E27").Value = "  -1.94%  "

$ws.Range("D28").Value = "19.39"
$ws.Range("E28").Value = "  -2.90%  "

$ws.Range("D29").Value = "6.004"
$ws.Range("E29").Value = "  -5.36%  "

$ws.Range("D30").Value = "2.005"
$ws.Range("E30").Value = "  -3.66%  "

$ws.Range("D31").Value = "117.58"
$ws.Range("E31").Value = "  -2.41%  "

$ws.Range("D32").Value = "0.9685"
$ws.Range("E32").Value = "  -3.28%  "

$ws.Range("D33").Value = "0.09394"
$ws.Range("E33").Value = "  -2.07%  "

$ws.Range("D34").Value = "5.373"
$ws.Range("E34").Value = "  -3.96%  "

$ws.Range("D35").Value = "3.463"
$ws.Range("E35").Value = "  -2.64%  "

$ws.Range("D36").Value = "1.343"
$ws.Range("E36").Value = "  -3.28%  "

$ws.Range("D37").Value = "0.06084"
$ws.Range("E37").Value = "  -6.43%  "

$ws.Range("D38").Value = "0.02191"
$ws.Range("E38").Value = "  -3.90%  "

$ws.Range("D39").Value = "1.164"
$ws.Range("E39").Value = "  -3.91%  "

$ws.Range("D40").Value = "0.5684"
$ws.Range("E40").Value = "  -3.94%  "

$ws.Range("E41").Value = "  -2.50%  "

$ws.Range("D42").Value = "10.14"
$ws.Range("E42").Value = "  -5.41%  "

$ws.Range("D43").Value = "0.1794"
$ws.Range("E43").Value = "  -2.55%  "

$ws.Range("D44").Value = "2.406"
$ws.Range("E44").Value = "  -4.31%  "

$ws.Range("D45").Value = "1.225"
$ws.Range("E45").Value = "  -4.48%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5367"
$ws.Range("E46").Value = "  -3.04%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "11.72"
$ws.Range("E47").Value = "  -4.96%  "

$ws.Range("D48").Value = "0.07096"
$ws.Range("E48").Value = "  -5.74%  "

$ws.Range("D49").Value = "1.900"
$ws.Range("E49").Value = "  -2.73%  "

$ws.Range("D50").Value = "113.66"
$ws.Range("E50").Value = "  -3.33%  "

$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -1.06%  "

# Restore default cell style on the Price column now that the
# text values are safely stored (matches the original workbook,
# which had no explicit number-format style on these cells).
$ws.Range("D2:D51").Style = "Normal"
